$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report header text (Volume/Number, and week-of date range) ---
$ws.Range("A8").Value = "Volume 30   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/16/2023  Through  1/22/2023"

# --- Style-switching cells: copy formatting+value from a stable donor cell, then set numeric value when needed ---
$ws.Range("J14").Copy($ws.Range("D15"))
$ws.Range("D15").Value2 = 1
$ws.Range("K14").Copy($ws.Range("E15"))
$ws.Range("E15").Value2 = 0
$ws.Range("I14").Copy($ws.Range("D22"))
$ws.Range("L14").Copy($ws.Range("E22"))
$ws.Range("I14").Copy($ws.Range("C23"))
$ws.Range("I14").Copy($ws.Range("D23"))
$ws.Range("L14").Copy($ws.Range("E23"))
$ws.Range("J14").Copy($ws.Range("D26"))
$ws.Range("D26").Value2 = 1
$ws.Range("K14").Copy($ws.Range("E26"))
$ws.Range("E26").Value2 = 100
$ws.Range("J14").Copy($ws.Range("D27"))
$ws.Range("D27").Value2 = 1
$ws.Range("K14").Copy($ws.Range("E27"))
$ws.Range("E27").Value2 = 0
$ws.Range("J14").Copy($ws.Range("C28"))
$ws.Range("C28").Value2 = 1
$ws.Range("J14").Copy($ws.Range("I28"))
$ws.Range("I28").Value2 = 1
$ws.Range("J14").Copy($ws.Range("C29"))
$ws.Range("C29").Value2 = 1
$ws.Range("J14").Copy($ws.Range("I29"))
$ws.Range("I29").Value2 = 1

# --- Plain numeric value updates (style already correct) ---
$ws.Range("F14").Value2 = 1
$ws.Range("G14").Value2 = 1
$ws.Range("H14").Value2 = 0
$ws.Range("C15").Value2 = 1
$ws.Range("F15").Value2 = 4
$ws.Range("G15").Value2 = 3
$ws.Range("H15").Value2 = 33.333333333333
$ws.Range("I15").Value2 = 4
$ws.Range("J15").Value2 = 3
$ws.Range("K15").Value2 = 33.333333333333
$ws.Range("L15").Value2 = 300
$ws.Range("M15").Value2 = 100
$ws.Range("N15").Value2 = 100
$ws.Range("C16").Value2 = 4
$ws.Range("D16").Value2 = 5
$ws.Range("E16").Value2 = -20
$ws.Range("F16").Value2 = 26
$ws.Range("G16").Value2 = 28
$ws.Range("H16").Value2 = -7.142857142857
$ws.Range("I16").Value2 = 19
$ws.Range("J16").Value2 = 23
$ws.Range("K16").Value2 = -17.391304347826
$ws.Range("L16").Value2 = 18.75
$ws.Range("M16").Value2 = -32.142857142857
$ws.Range("N16").Value2 = -84.552845528455
$ws.Range("C17").Value2 = 15
$ws.Range("D17").Value2 = 3
$ws.Range("E17").Value2 = 400
$ws.Range("F17").Value2 = 45
$ws.Range("G17").Value2 = 32
$ws.Range("H17").Value2 = 40.625
$ws.Range("I17").Value2 = 34
$ws.Range("J17").Value2 = 27
$ws.Range("K17").Value2 = 25.925925925925
$ws.Range("L17").Value2 = 54.545454545454
$ws.Range("M17").Value2 = 161.538461538462
$ws.Range("N17").Value2 = -37.037037037037
$ws.Range("C18").Value2 = 2
$ws.Range("D18").Value2 = 4
$ws.Range("E18").Value2 = -50
$ws.Range("F18").Value2 = 15
$ws.Range("G18").Value2 = 17
$ws.Range("H18").Value2 = -11.764705882352
$ws.Range("I18").Value2 = 13
$ws.Range("J18").Value2 = 12
$ws.Range("K18").Value2 = 8.333333333333
$ws.Range("L18").Value2 = 30
$ws.Range("M18").Value2 = 18.181818181818
$ws.Range("N18").Value2 = -80.30303030303
$ws.Range("C19").Value2 = 7
$ws.Range("D19").Value2 = 10
$ws.Range("E19").Value2 = -30
$ws.Range("F19").Value2 = 36
$ws.Range("G19").Value2 = 56
$ws.Range("H19").Value2 = -35.714285714285
$ws.Range("I19").Value2 = 30
$ws.Range("J19").Value2 = 40
$ws.Range("K19").Value2 = -25
$ws.Range("L19").Value2 = 66.666666666666
$ws.Range("M19").Value2 = 20
$ws.Range("N19").Value2 = -52.380952380952
$ws.Range("C20").Value2 = 4
$ws.Range("D20").Value2 = 8
$ws.Range("E20").Value2 = -50
$ws.Range("F20").Value2 = 19
$ws.Range("G20").Value2 = 27
$ws.Range("H20").Value2 = -29.629629629629
$ws.Range("I20").Value2 = 15
$ws.Range("J20").Value2 = 26
$ws.Range("K20").Value2 = -42.307692307692
$ws.Range("L20").Value2 = 114.285714285714
$ws.Range("M20").Value2 = 114.285714285714
$ws.Range("N20").Value2 = -82.142857142857
$ws.Range("C21").Value2 = 33
$ws.Range("D21").Value2 = 31
$ws.Range("E21").Value2 = 6.451612903225
$ws.Range("F21").Value2 = 146
$ws.Range("G21").Value2 = 164
$ws.Range("H21").Value2 = -10.975609756097
$ws.Range("I21").Value2 = 115
$ws.Range("J21").Value2 = 132
$ws.Range("K21").Value2 = -12.878787878787
$ws.Range("L21").Value2 = 55.405405405405
$ws.Range("M21").Value2 = 33.720930232558
$ws.Range("N21").Value2 = -70.663265306122
$ws.Range("F23").Value2 = 2
$ws.Range("G23").Value2 = 5
$ws.Range("H23").Value2 = -60
$ws.Range("M23").Value2 = 0
$ws.Range("C24").Value2 = 29
$ws.Range("D24").Value2 = 30
$ws.Range("E24").Value2 = -3.333333333333
$ws.Range("F24").Value2 = 111
$ws.Range("G24").Value2 = 116
$ws.Range("H24").Value2 = -4.310344827586
$ws.Range("I24").Value2 = 89
$ws.Range("J24").Value2 = 93
$ws.Range("K24").Value2 = -4.301075268817
$ws.Range("L24").Value2 = 41.269841269841
$ws.Range("M24").Value2 = 36.923076923076
$ws.Range("C25").Value2 = 13
$ws.Range("D25").Value2 = 7
$ws.Range("E25").Value2 = 85.714285714285
$ws.Range("F25").Value2 = 72
$ws.Range("G25").Value2 = 39
$ws.Range("H25").Value2 = 84.615384615384
$ws.Range("I25").Value2 = 57
$ws.Range("J25").Value2 = 32
$ws.Range("K25").Value2 = 78.125
$ws.Range("L25").Value2 = 137.5
$ws.Range("M25").Value2 = 5.555555555555
$ws.Range("F26").Value2 = 5
$ws.Range("H26").Value2 = 25
$ws.Range("I26").Value2 = 5
$ws.Range("J26").Value2 = 4
$ws.Range("K26").Value2 = 25
$ws.Range("L26").Value2 = 400
$ws.Range("C27").Value2 = 1
$ws.Range("G27").Value2 = 3
$ws.Range("H27").Value2 = 233.333333333333
$ws.Range("I27").Value2 = 5
$ws.Range("J27").Value2 = 3
$ws.Range("K27").Value2 = 66.666666666666
$ws.Range("L27").Value2 = 25
$ws.Range("F28").Value2 = 2
$ws.Range("G28").Value2 = 2
$ws.Range("K28").Value2 = 0
$ws.Range("N28").Value2 = -90
$ws.Range("F29").Value2 = 2
$ws.Range("G29").Value2 = 2
$ws.Range("K29").Value2 = 0
$ws.Range("N29").Value2 = -88.888888888888